$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 3005033
$ws.Cells.Item(34, 9).Value = 3600639.5
$ws.Cells.Item(34, 10).Value = 27000
$ws.Cells.Item(34, 11).Value = 3600639.5
$ws.Cells.Item(34, 12).Value = 27000
$ws.Cells.Item(34, 13).Value = -3600436.5
$ws.Cells.Item(34, 14).Value = -27406
$ws.Cells.Item(36, 8).Value = 3005033
$ws.Cells.Item(36, 9).Value = 3600639.5
$ws.Cells.Item(36, 10).Value = 27000
$ws.Cells.Item(36, 11).Value = 3600639.5
$ws.Cells.Item(36, 12).Value = 27000
$ws.Cells.Item(36, 13).Value = -3599924.5
$ws.Cells.Item(36, 14).Value = -28430
$ws.Cells.Item(112, 8).Value = 1407.4286
$ws.Cells.Item(112, 10).Value = 1850.4
$ws.Cells.Item(112, 12).Value = 5551.200000000001
$ws.Cells.Item(112, 14).Value = -7767.200000000001
$ws.Cells.Item(137, 8).Value = 1952.6216
$ws.Cells.Item(137, 9).Value = 1663.9688
$ws.Cells.Item(137, 11).Value = 4991.9064
$ws.Cells.Item(137, 13).Value = -2441.9064
$ws.Cells.Item(141, 8).Value = 4246.263
$ws.Cells.Item(141, 9).Value = 1619.3334
$ws.Cells.Item(141, 10).Value = 6610.5
$ws.Cells.Item(141, 11).Value = 4858.0002
$ws.Cells.Item(141, 12).Value = 19831.5
$ws.Cells.Item(141, 13).Value = 321.9997999999996
$ws.Cells.Item(141, 14).Value = -30191.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 10085.303
$ws.Cells.Item(61, 9).Value = 5807.3184
$ws.Cells.Item(61, 10).Value = 18641.273
$ws.Cells.Item(61, 11).Value = 5807.3184
$ws.Cells.Item(61, 12).Value = 18641.273
$ws.Cells.Item(61, 13).Value = -5595.3184
$ws.Cells.Item(61, 14).Value = -19065.273
$ws.Cells.Item(74, 8).Value = 4996.353
$ws.Cells.Item(74, 9).Value = 2093.6128
$ws.Cells.Item(74, 10).Value = 34991.332
$ws.Cells.Item(74, 11).Value = 2093.6128
$ws.Cells.Item(74, 12).Value = 34991.332
$ws.Cells.Item(74, 13).Value = -1219.6128
$ws.Cells.Item(74, 14).Value = -36739.332
$ws.Cells.Item(77, 8).Value = 4996.353
$ws.Cells.Item(77, 9).Value = 2093.6128
$ws.Cells.Item(77, 10).Value = 34991.332
$ws.Cells.Item(77, 11).Value = 10468.064
$ws.Cells.Item(77, 12).Value = 174956.66
$ws.Cells.Item(77, 13).Value = -6100.063999999998
$ws.Cells.Item(77, 14).Value = -183692.66
$ws.Cells.Item(124, 8).Value = 26143
$ws.Cells.Item(124, 10).Value = 26143
$ws.Cells.Item(124, 12).Value = 26143
$ws.Cells.Item(124, 14).Value = -35963
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2557.375
$ws.Cells.Item(132, 9).Value = 1668.9286
$ws.Cells.Item(132, 11).Value = 5006.7858
$ws.Cells.Item(132, 13).Value = -2476.7858
$ws.Cells.Item(136, 8).Value = 10085.303
$ws.Cells.Item(136, 9).Value = 5807.3184
$ws.Cells.Item(136, 10).Value = 18641.273
$ws.Cells.Item(136, 11).Value = 17421.9552
$ws.Cells.Item(136, 12).Value = 55923.819
$ws.Cells.Item(136, 13).Value = -14871.9552
$ws.Cells.Item(136, 14).Value = -61023.819
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3746.818
$ws.Cells.Item(20, 9).Value = 3652
$ws.Cells.Item(20, 10).Value = 3999.6667
$ws.Cells.Item(20, 11).Value = 3652
$ws.Cells.Item(20, 12).Value = 3999.6667
$ws.Cells.Item(20, 13).Value = -3405
$ws.Cells.Item(20, 14).Value = -4493.6667
$ws.Cells.Item(134, 8).Value = 34525.805
$ws.Cells.Item(134, 9).Value = 2453.2222
$ws.Cells.Item(134, 11).Value = 7359.6666
$ws.Cells.Item(134, 13).Value = -4824.6666
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6349.185
$ws.Cells.Item(31, 9).Value = 8579
$ws.Cells.Item(31, 11).Value = 8579
$ws.Cells.Item(31, 13).Value = -8284
$ws.Cells.Item(34, 8).Value = 6349.185
$ws.Cells.Item(34, 9).Value = 8579
$ws.Cells.Item(34, 11).Value = 8579
$ws.Cells.Item(34, 13).Value = -8377
$ws.Cells.Item(58, 8).Value = 3954814.5
$ws.Cells.Item(58, 9).Value = 5683397.5
$ws.Cells.Item(58, 10).Value = 3767.7144
$ws.Cells.Item(58, 11).Value = 5683397.5
$ws.Cells.Item(58, 12).Value = 3767.7144
$ws.Cells.Item(58, 13).Value = -5683194.5
$ws.Cells.Item(58, 14).Value = -4173.7144
$ws.Cells.Item(63, 8).Value = 32895.855
$ws.Cells.Item(63, 10).Value = 32895.855
$ws.Cells.Item(63, 12).Value = 32895.855
$ws.Cells.Item(63, 14).Value = -34267.855
$ws.Cells.Item(66, 8).Value = 32895.855
$ws.Cells.Item(66, 10).Value = 32895.855
$ws.Cells.Item(66, 12).Value = 98687.565
$ws.Cells.Item(66, 14).Value = -105551.565
$ws.Cells.Item(132, 8).Value = 2939.3142
$ws.Cells.Item(132, 9).Value = 2559.12
$ws.Cells.Item(132, 10).Value = 3889.8
$ws.Cells.Item(132, 11).Value = 7677.36
$ws.Cells.Item(132, 12).Value = 11669.4
$ws.Cells.Item(132, 13).Value = -5147.36
$ws.Cells.Item(132, 14).Value = -16729.4
$ws.Cells.Item(134, 8).Value = 2603.0605
$ws.Cells.Item(134, 9).Value = 2598.0476
$ws.Cells.Item(134, 10).Value = 2611.8333
$ws.Cells.Item(134, 11).Value = 7794.1428
$ws.Cells.Item(134, 12).Value = 7835.499899999999
$ws.Cells.Item(134, 13).Value = -5259.1428
$ws.Cells.Item(134, 14).Value = -12905.4999
$ws.Cells.Item(136, 8).Value = 3954814.5
$ws.Cells.Item(136, 9).Value = 5683397.5
$ws.Cells.Item(136, 10).Value = 3767.7144
$ws.Cells.Item(136, 11).Value = 17050192.5
$ws.Cells.Item(136, 12).Value = 11303.1432
$ws.Cells.Item(136, 13).Value = -17047642.5
$ws.Cells.Item(136, 14).Value = -16403.1432
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 19750
$ws.Cells.Item(5, 10).Value = 19750
$ws.Cells.Item(5, 12).Value = 19750
$ws.Cells.Item(5, 14).Value = -19974
$ws.Cells.Item(33, 8).Value = 9800
$ws.Cells.Item(33, 10).Value = 9800
$ws.Cells.Item(33, 12).Value = 9800
$ws.Cells.Item(33, 14).Value = -10304
$ws.Cells.Item(70, 8).Value = 5671.66
$ws.Cells.Item(70, 10).Value = 5795.2383
$ws.Cells.Item(70, 12).Value = 5795.2383
$ws.Cells.Item(70, 14).Value = -6335.2383
$ws.Cells.Item(73, 8).Value = 5671.66
$ws.Cells.Item(73, 10).Value = 5795.2383
$ws.Cells.Item(73, 12).Value = 5795.2383
$ws.Cells.Item(73, 14).Value = -7667.2383
$ws.Cells.Item(126, 8).Value = 2863.2632
$ws.Cells.Item(126, 9).Value = 1822.4445
$ws.Cells.Item(126, 11).Value = 5467.333500000001
$ws.Cells.Item(126, 13).Value = -2997.333500000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 341428.34
$ws.Cells.Item(61, 9).Value = 12117.895
$ws.Cells.Item(61, 10).Value = 910237.25
$ws.Cells.Item(61, 11).Value = 12117.895
$ws.Cells.Item(61, 12).Value = 910237.25
$ws.Cells.Item(61, 13).Value = -11915.895
$ws.Cells.Item(61, 14).Value = -910641.25
$ws.Cells.Item(113, 8).Value = 341428.34
$ws.Cells.Item(113, 9).Value = 12117.895
$ws.Cells.Item(113, 10).Value = 910237.25
$ws.Cells.Item(113, 11).Value = 12117.895
$ws.Cells.Item(113, 12).Value = 910237.25
$ws.Cells.Item(113, 13).Value = -9947.895
$ws.Cells.Item(113, 14).Value = -914577.25
$ws.Cells.Item(132, 8).Value = 3667.9143
$ws.Cells.Item(132, 9).Value = 3439.1924
$ws.Cells.Item(132, 10).Value = 4328.6665
$ws.Cells.Item(132, 11).Value = 10317.5772
$ws.Cells.Item(132, 12).Value = 12985.9995
$ws.Cells.Item(132, 13).Value = -7787.5772
$ws.Cells.Item(132, 14).Value = -18045.9995
$ws.Cells.Item(136, 8).Value = 6049.875
$ws.Cells.Item(136, 9).Value = 3969.1333
$ws.Cells.Item(136, 10).Value = 7885.8237
$ws.Cells.Item(136, 11).Value = 11907.3999
$ws.Cells.Item(136, 12).Value = 23657.4711
$ws.Cells.Item(136, 13).Value = -9357.3999
$ws.Cells.Item(136, 14).Value = -28757.4711
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 12034.714
$ws.Cells.Item(2, 9).Value = 1268.5
$ws.Cells.Item(2, 10).Value = 16341.2
$ws.Cells.Item(2, 11).Value = 1268.5
$ws.Cells.Item(2, 12).Value = 16341.2
$ws.Cells.Item(2, 13).Value = -1156.5
$ws.Cells.Item(2, 14).Value = -16565.2
$ws.Cells.Item(53, 8).Value = 10000
$ws.Cells.Item(53, 10).Value = 10000
$ws.Cells.Item(53, 12).Value = 10000
$ws.Cells.Item(53, 14).Value = -11214
$ws.Cells.Item(132, 8).Value = 1907.9231
$ws.Cells.Item(132, 9).Value = 1134.4615
$ws.Cells.Item(132, 11).Value = 3403.3845
$ws.Cells.Item(132, 13).Value = -873.3844999999997
$ws.Cells.Item(136, 8).Value = 6722.2188
$ws.Cells.Item(136, 9).Value = 2822.5557
$ws.Cells.Item(136, 10).Value = 11736.071
$ws.Cells.Item(136, 11).Value = 8467.667099999999
$ws.Cells.Item(136, 12).Value = 35208.213
$ws.Cells.Item(136, 13).Value = -5917.667099999999
$ws.Cells.Item(136, 14).Value = -40308.213
